# Moved new hardcoded dialogue into the spreadsheet
# Adds 3 new rows (13-15) to the "dialogue" sheet for:
#   - deleteDataButton / Delete Data button text
#   - instagram link label
#   - github link label
# Columns: A=key, B=english, C=czech, D=korean, E=japanese

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Delete Data button
$ws.Range("A13").Value = "deleteDataButton"
$ws.Range("B13").Value = "Delete Data"
$ws.Range("C13").Value = "Smazat Data"
$ws.Range("D13").Value = "데이터 삭제"
$ws.Range("E13").Value = "データ削除"

# Row 14: Instagram link
$ws.Range("A14").Value = "instagram"
$ws.Range("B14").Value = "My Instagram"
$ws.Range("C14").Value = "Můj Instagram"
$ws.Range("D14").Value = "내 인스타그램"
$ws.Range("E14").Value = "私のインスタグラム"

# Row 15: GitHub link
$ws.Range("A15").Value = "github"
$ws.Range("B15").Value = "Project GitHub"
$ws.Range("C15").Value = "GitHub Projektu"
$ws.Range("D15").Value = "프로젝트 깃허브"
$ws.Range("E15").Value = "プロジェクト ギツハブ"

# Match the row height used by the rest of the table
$ws.Rows.Item(13).RowHeight = 42
$ws.Rows.Item(14).RowHeight = 42
$ws.Rows.Item(15).RowHeight = 42

# Move the selection/scroll to reflect where editing ended up
[void]$ws.Range("C17").Select()
